$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (rich-text / multi-run shared strings).
#    Only the specific run's text changes; surrounding runs / formatting
#    are left alone by targeting the exact sub-string with .Characters().
# ---------------------------------------------------------------------------

# A8 holds: "Volume " + "32" + "   Number  " + "5"  -> bump issue number 5 -> 6
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$numStart = $volText.IndexOf("   Number  ") + "   Number  ".Length + 1
$numLen = $volText.Length - ($numStart - 1)
$volCell.Characters($numStart, $numLen).Text = "6"

# C9 holds: "Report Covering the Week  " + "1/27/2025" + "  Through  " + "2/2/2025"
# -> shift the reporting week forward by one week.
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value2
$d1Start = $weekText.IndexOf("1/27/2025") + 1
$weekCell.Characters($d1Start, "1/27/2025".Length).Text = "2/3/2025"

$weekText2 = $weekCell.Value2
$d2Start = $weekText2.IndexOf("2/2/2025") + 1
$weekCell.Characters($d2Start, "2/2/2025".Length).Text = "2/9/2025"

# ---------------------------------------------------------------------------
# 2) A handful of cells in rows 31 & 33 flip from the "no data" placeholder
#    (text "0" / "***.*") to real numeric figures now that data exists for
#    those precincts. Give them the same numeric styles used elsewhere in
#    the table (count column / pct-change column) before writing values.
# ---------------------------------------------------------------------------

$countFormat = "#,##0"
$pctFormat = '#,##0.0;"-"#,##0.0'

$countCells = @("D31", "G31", "C33", "D33", "F33", "I33")
foreach ($cell in $countCells) {
    $ws.Range($cell).NumberFormat = $countFormat
}

$pctCells = @("E31", "H31", "E33")
foreach ($cell in $pctCells) {
    $ws.Range($cell).NumberFormat = $pctFormat
}

# ---------------------------------------------------------------------------
# 3) Refreshed weekly crime-statistics figures (counts + computed % changes)
#    for rows 14-30 and 31/33, per this week's CompStat data pull.
# ---------------------------------------------------------------------------

$values = @(
    @("D14", 3),
    @("E14", -66.666666666666),
    @("F14", 3),
    @("G14", 8),
    @("H14", -62.5),
    @("I14", 9),
    @("J14", 12),
    @("K14", -25),
    @("L14", -25),
    @("M14", 12.5),
    @("N14", -82.692307692307),
    @("C15", 8),
    @("D15", 9),
    @("E15", -11.111111111111),
    @("F15", 32),
    @("G15", 34),
    @("H15", -5.882352941176),
    @("I15", 47),
    @("J15", 45),
    @("K15", 4.444444444444),
    @("L15", 4.444444444444),
    @("M15", 46.875),
    @("N15", -18.965517241379),
    @("C16", 74),
    @("D16", 85),
    @("E16", -12.941176470588),
    @("F16", 309),
    @("G16", 377),
    @("H16", -18.037135278514),
    @("I16", 440),
    @("J16", 560),
    @("K16", -21.428571428571),
    @("L16", -10.020449897750),
    @("M16", -4.761904761904),
    @("N16", -76.470588235294),
    @("C17", 142),
    @("D17", 135),
    @("E17", 5.185185185185),
    @("F17", 566),
    @("G17", 538),
    @("H17", 5.204460966542),
    @("I17", 790),
    @("J17", 779),
    @("K17", 1.412066752246),
    @("L17", 7.629427792915),
    @("M17", 85.882352941176),
    @("N17", -1.618929016189),
    @("C18", 50),
    @("D18", 52),
    @("E18", -3.846153846153),
    @("F18", 215),
    @("G18", 216),
    @("H18", -0.462962962962),
    @("I18", 296),
    @("J18", 299),
    @("K18", -1.003344481605),
    @("L18", -12.166172106824),
    @("M18", -21.899736147757),
    @("N18", -85.631067961165),
    @("C19", 160),
    @("D19", 187),
    @("E19", -14.438502673796),
    @("F19", 642),
    @("G19", 757),
    @("H19", -15.191545574636),
    @("I19", 866),
    @("J19", 1043),
    @("K19", -16.970278044103),
    @("L19", 11.454311454311),
    @("M19", 100.462962962963),
    @("N19", 16.711590296496),
    @("C20", 78),
    @("D20", 69),
    @("E20", 13.043478260869),
    @("F20", 308),
    @("G20", 316),
    @("H20", -2.531645569620),
    @("I20", 426),
    @("J20", 449),
    @("K20", -5.122494432071),
    @("L20", -30.956239870340),
    @("M20", 102.857142857143),
    @("N20", -75.275681950087),
    @("C21", 513),
    @("D21", 540),
    @("E21", -5),
    @("F21", 2075),
    @("G21", 2246),
    @("H21", -7.613535173642),
    @("I21", 2874),
    @("J21", 3187),
    @("K21", -9.821148415437),
    @("L21", -4.549983394221),
    @("M21", 47.535934291581),
    @("N21", -60.673234811165),
    @("C22", 4),
    @("D22", 3),
    @("E22", 33.333333333333),
    @("G22", 28),
    @("H22", -7.142857142857),
    @("I22", 34),
    @("J22", 41),
    @("K22", -17.073170731707),
    @("L22", 54.545454545454),
    @("M22", 9.677419354838),
    @("C23", 20),
    @("D23", 30),
    @("E23", -33.333333333333),
    @("F23", 99),
    @("H23", -22.047244094488),
    @("I23", 145),
    @("J23", 187),
    @("K23", -22.459893048128),
    @("L23", -25.641025641025),
    @("M23", 40.776699029126),
    @("C24", 277),
    @("D24", 294),
    @("E24", -5.782312925170),
    @("F24", 1234),
    @("G24", 1304),
    @("H24", -5.368098159509),
    @("I24", 1739),
    @("J24", 1823),
    @("K24", -4.607789358200),
    @("L24", -1.806888763410),
    @("M24", 31.842304776345),
    @("C25", 88),
    @("D25", 117),
    @("E25", -24.786324786324),
    @("F25", 421),
    @("G25", 559),
    @("H25", -24.686940966010),
    @("I25", 590),
    @("J25", 803),
    @("K25", -26.525529265255),
    @("L25", -29.086538461538),
    @("C26", 182),
    @("D26", 191),
    @("E26", -4.712041884816),
    @("F26", 717),
    @("G26", 768),
    @("H26", -6.640625),
    @("I26", 997),
    @("J26", 1083),
    @("K26", -7.940904893813),
    @("L26", -1.287128712871),
    @("M26", -0.3),
    @("D27", 12),
    @("E27", 0),
    @("F27", 42),
    @("G27", 54),
    @("H27", -22.222222222222),
    @("I27", 62),
    @("J27", 72),
    @("K27", -13.888888888888),
    @("L27", -15.068493150684),
    @("C28", 13),
    @("D28", 25),
    @("E28", -48),
    @("F28", 79),
    @("G28", 77),
    @("H28", 2.597402597402),
    @("I28", 111),
    @("J28", 107),
    @("K28", 3.738317757009),
    @("L28", -6.722689075630),
    @("C29", 2),
    @("D29", 5),
    @("E29", -60),
    @("F29", 12),
    @("G29", 23),
    @("H29", -47.826086956521),
    @("I29", 24),
    @("J29", 34),
    @("K29", -29.411764705882),
    @("L29", -22.580645161290),
    @("M29", -38.461538461538),
    @("N29", -79.487179487179),
    @("C30", 2),
    @("D30", 5),
    @("E30", -60),
    @("F30", 11),
    @("H30", -45),
    @("I30", 20),
    @("J30", 29),
    @("K30", -31.034482758620),
    @("L30", -20),
    @("M30", -44.444444444444),
    @("N30", -81.481481481481),
    @("D31", 1),
    @("E31", -100),
    @("G31", 1),
    @("H31", -100),
    @("J31", 2),
    @("C33", 2),
    @("D33", 1),
    @("E33", 100),
    @("F33", 2),
    @("G33", 4),
    @("H33", -50),
    @("I33", 2),
    @("J33", 5),
    @("K33", -60),
    @("L33", 0)
)

foreach ($pair in $values) {
    $ws.Range($pair[0]).Value = $pair[1]
}
